$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 81

# Plain text values - these are not pure numeric strings so Excel keeps them as text.
$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# "25" looks numeric, so a bare assignment would store it as a number.
# Prefix with an apostrophe (quote-prefix) to force Excel to keep it as text,
# matching the source data which stores it as a string value.
$ws.Cells.Item($row, 3).Value = "'25"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
